# Generate Report for Handoff
# Update the localization-status workbook so that file "b.md" is reported
# as handed off (instead of "handed back"), with a new handoff xliff file,
# a new handoff timestamp, and an error detail describing the stale
# handback version - mirroring the change made for the real handoff run.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2565f4f8b28b1fb60927bb690b4be35636bb3f04/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/010894195651fb9c61622bc4d196526531d257d4/e2e/b.md."

# ---------------------------------------------------------------------
# Overview sheet: row 3 is b.md - status + latest handoff date refreshed
# ---------------------------------------------------------------------
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-15 12:31:47"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 is b.md
# ---------------------------------------------------------------------
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-15 12:31:42"
$zhcn.Range("P3").Value = $errorDetail

# ---------------------------------------------------------------------
# de-de sheet: row 3 is b.md
# ---------------------------------------------------------------------
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-15 12:31:47"
$dede.Range("P3").Value = $errorDetail

# ---------------------------------------------------------------------
# Widen column P ("Error Detail") on zh-cn / de-de now that it carries
# the long error message text.
# ---------------------------------------------------------------------
$zhcn.Range("P1").ColumnWidth = 40
$dede.Range("P1").ColumnWidth = 40
